$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 6; this shifts the former rows 6-41
# down to 7-42, preserving all of their existing data/styles.
$ws.Rows.Item(6).Insert()

# Columns that stay identical to the rest of the table (copy from the
# row directly below, which holds the data that used to be row 6).
$ws.Range("A6").Value = $ws.Range("A7").Value2
$ws.Range("B6").Value = $ws.Range("B7").Value2
$ws.Range("C6").Value = $ws.Range("C7").Value2
$ws.Range("E6").Value = $ws.Range("E7").Value2
$ws.Range("F6").Value = $ws.Range("F7").Value2
$ws.Range("G6").Value = $ws.Range("G7").Value2
$ws.Range("N6").Value = $ws.Range("N7").Value2
$ws.Range("O6").Value = $ws.Range("O7").Value2
$ws.Range("Q6").Value = $ws.Range("Q7").Value2
$ws.Range("R6").Value = $ws.Range("R7").Value2

# New values for the inserted row.
$ws.Range("D6").Value = 45163
$ws.Range("H6").Value = "Argentina(o)"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("P6").Value = 280
